$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "retention_days" header in F1, reusing the same header
# formatting (bold/border/centered) already applied to the other headers
# by copying E1's format onto F1.
$ws.Cells.Item(1, 6).Value = "retention_days"
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)

# Per-row updates: refreshed last_backup date (column C, kept as literal
# text via a leading apostrophe so it is not auto-converted to a date
# serial) and a new retention_days value (column F).
$rows = @(
    @{ Row = 2;  Date = "2025-12-01"; Retention = 30 },
    @{ Row = 3;  Date = "2025-12-10"; Retention = 30 },
    @{ Row = 4;  Date = "2025-11-29"; Retention = 30 },
    @{ Row = 5;  Date = "2025-12-08"; Retention = 30 },
    @{ Row = 6;  Date = "2025-11-27"; Retention = 15 },
    @{ Row = 7;  Date = "2025-12-11"; Retention = 30 },
    @{ Row = 8;  Date = "2025-11-25"; Retention = 30 },
    @{ Row = 9;  Date = "2025-12-09"; Retention = 30 },
    @{ Row = 10; Date = "2025-11-23"; Retention = 30 },
    @{ Row = 11; Date = "2025-12-12"; Retention = 15 },
    @{ Row = 12; Date = "2025-11-21"; Retention = 30 },
    @{ Row = 13; Date = "2025-12-10"; Retention = 30 },
    @{ Row = 14; Date = "2025-11-19"; Retention = 30 },
    @{ Row = 15; Date = "2025-12-08"; Retention = 30 },
    @{ Row = 16; Date = "2025-11-17"; Retention = 15 },
    @{ Row = 17; Date = "2025-12-11"; Retention = 30 },
    @{ Row = 18; Date = "2025-11-15"; Retention = 30 },
    @{ Row = 19; Date = "2025-12-09"; Retention = 30 },
    @{ Row = 20; Date = "2025-11-13"; Retention = 30 },
    @{ Row = 21; Date = "2025-12-12"; Retention = 15 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = "'" + $r.Date
    $ws.Cells.Item($r.Row, 6).Value = $r.Retention
}
